$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# =========================================================================
# Auto data update: refreshed win-rate percentages scraped from MAA copilot
# operator list, plus a few newly catalogued operators (CI commit #75).
# =========================================================================

# --- 1) In-place refresh of existing "maa://..." win-rate strings ---
$ws.Range("AE2").Value = 'maa://25251 (92.41), ***maa://21730 (17.19), ***maa://39501 (18.18), *maa://36675 (60.0)'
$ws.Range("K3").Value = '*maa://22880 (70.2), maa://20276 (82.73), *maa://22749 (62.5)'
$ws.Range("O3").Value = 'maa://21249 (95.12), maa://26254 (95.24)'
$ws.Range("S3").Value = 'maa://24617 (88.35), **maa://20790 (43.94), ***maa://37170 (20.0)'
$ws.Range("W3").Value = 'maa://27396 (84.81), maa://27484 (95.74), maa://27480 (82.35)'
$ws.Range("W4").Value = '**maa://32495 (47.93), ***maa://31785 (16.51), ***maa://36683 (26.67)'
$ws.Range("AE4").Value = '*maa://30062 (61.9), ***maa://26209 (13.04), *maa://39394 (78.57)'
$ws.Range("AA5").Value = '*maa://29863 (75.0), ***maa://22752 (13.33), **maa://26013 (42.86)'
$ws.Range("G6").Value = 'maa://24370 (96.15)'
$ws.Range("G7").Value = '*maa://22763 (66.67)'
$ws.Range("W7").Value = 'maa://22399 (94.62), *maa://22758 (70.59)'
$ws.Range("AE7").Value = '*maa://26191 (68.49), *maa://36671 (72.09)'
$ws.Range("C8").Value = '*maa://21476 (69.77), ***maa://39431 (25.0), *maa://37551 (57.14)'
$ws.Range("W8").Value = 'maa://21411 (96.31)'
$ws.Range("AE9").Value = 'maa://26206 (88.89), **maa://22865 (45.65)'
$ws.Range("C10").Value = '***maa://25695 (19.3), **maa://32237 (37.84), ***maa://34206 (14.29), ***maa://39951 (21.74), **maa://39243 (40.0)'
$ws.Range("S10").Value = 'maa://27395 (96.64), maa://22755 (87.62), **maa://22756 (40.91), ***maa://21737 (10.61)'
$ws.Range("W10").Value = 'maa://22301 (97.38), maa://22726 (100.0)'
$ws.Range("K11").Value = 'maa://21287 (87.36)'
$ws.Range("W11").Value = 'maa://36713 (97.83)'
$ws.Range("G12").Value = 'maa://21867 (89.93)'
$ws.Range("W12").Value = 'maa://22753 (91.84), *maa://21485 (76.56), maa://37962 (81.25)'
$ws.Range("AA12").Value = 'maa://23669 (95.86), maa://36677 (94.87), maa://39872 (83.33)'
$ws.Range("C13").Value = 'maa://24999 (91.41), maa://36673 (91.8), maa://25001 (85.51)'
$ws.Range("W13").Value = '*maa://34957 (75.0), *maa://22768 (53.33)'
$ws.Range("C14").Value = 'maa://30764 (85.71)'
$ws.Range("S14").Value = '*maa://22471 (59.42), maa://22521 (94.44)'
$ws.Range("C15").Value = '*maa://22743 (76.19), maa://22734 (83.33), *maa://30808 (63.64), ***maa://36048 (13.33)'
$ws.Range("C16").Value = 'maa://21441 (96.17), maa://36679 (93.55), maa://37650 (95.45)'
$ws.Range("S16").Value = 'maa://22729 (95.17), *maa://28648 (69.09), *maa://36674 (78.57)'
$ws.Range("AE16").Value = '*maa://23911 (61.54), maa://27755 (91.78)'
$ws.Range("C18").Value = 'maa://24570 (96.49)'
$ws.Range("K18").Value = 'maa://22466 (88.37), *maa://22732 (52.5)'
$ws.Range("S19").Value = 'maa://24386 (98.75)'
$ws.Range("K20").Value = 'maa://41331 (90.32)'
$ws.Range("AE21").Value = 'maa://22524 (94.29), *maa://22432 (74.07)'
$ws.Range("W22").Value = 'maa://21282 (98.82), *maa://37649 (64.71)'
$ws.Range("K23").Value = 'maa://39756 (92.11), maa://39875 (95.56)'
$ws.Range("W24").Value = 'maa://23504 (92.86), maa://29988 (86.0), **maa://22892 (40.14), *maa://25141 (76.86), *maa://36663 (80.0), ***maa://22815 (23.08)'
$ws.Range("AE24").Value = 'maa://22523 (84.86), *maa://36672 (76.74), maa://29910 (93.88), **maa://21440 (34.55)'
$ws.Range("K25").Value = 'maa://24378 (88.89)'
$ws.Range("O25").Value = 'maa://24382 (92.0)'
$ws.Range("G26").Value = 'maa://24913 (91.04)'
$ws.Range("C28").Value = 'maa://24465 (90.32), maa://25725 (82.28)'
$ws.Range("S28").Value = 'maa://23263 (94.25), *maa://29765 (62.12)'
$ws.Range("W28").Value = 'maa://39929 (86.02), ***maa://39723 (14.71), maa://41749 (85.71)'
$ws.Range("AE28").Value = 'maa://36660 (93.75), *maa://36701 (64.0)'
$ws.Range("K29").Value = 'maa://28432 (93.47), *maa://28440 (72.84), maa://31400 (100.0), *maa://28650 (66.67)'
$ws.Range("O29").Value = '*maa://23168 (55.77), **maa://30050 (40.0)'
$ws.Range("K31").Value = 'maa://35926 (93.53), maa://36258 (80.52)'
$ws.Range("G32").Value = 'maa://21895 (97.01), maa://36667 (98.04), **maa://20793 (38.78), maa://22760 (100.0)'
$ws.Range("S32").Value = 'maa://41108 (88.89), maa://41238 (94.12)'
$ws.Range("K35").Value = 'maa://41296 (97.5)'
$ws.Range("AE35").Value = 'maa://39479 (90.0)'
$ws.Range("AE38").Value = 'maa://36697 (85.71)'
$ws.Range("O40").Value = 'maa://23278 (96.21), maa://21386 (95.63), maa://36664 (92.31)'
$ws.Range("G44").Value = 'maa://29768 (97.51), maa://27728 (96.0)'
$ws.Range("S44").Value = 'maa://39366 (86.36)'
$ws.Range("S45").Value = '*maa://39364 (62.5)'
$ws.Range("G60").Value = '**maa://40438 (37.5)'

# --- 2) "-" placeholders that become literal "0" / "None" text.
#        A leading apostrophe forces Excel to store these as TEXT instead
#        of re-parsing "0" as the number 0.
$ws.Range("R17").Value = "'0"
$ws.Range("S17").Value = 'None'
$ws.Range("Z26").Value = "'0"
$ws.Range("AA26").Value = 'None'
$ws.Range("J36").Value = "'0"
$ws.Range("K36").Value = 'None'

# --- 3) Newly added operators in the "辅助" (Support) block (rows 36-37) ---
$ws.Range("Y36").Value = "凯瑟琳"
$ws.Range("Z36").Value = "-"
$ws.Range("AA36").Value = "-"

$ws.Range("Y37").Value = "波卜"
$ws.Range("Z37").Value = "-"
$ws.Range("AA37").Value = "-"

# --- 4) Newly added operator row 73 ("近卫"/Guard block: 维娜·维多利亚) ---
$ws.Range("E73").Value = "维娜·维多利亚"
$ws.Range("F73").Value = "-"
$ws.Range("G73").Value = "-"
